# Update daily and weekly charts
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Regional Bed Avaliability")
$ws2 = $wb.Worksheets.Item("Hospital COVID Census")

# ---- Sheet 1: "Regional Bed Avaliability" ----
$ws1.Range("C2").Value = 2700
$ws1.Range("E2").Value = 100
$ws1.Range("F2").Value = 500

$ws1.Range("F5").Value = 200

$ws1.Range("C7").Value = 900
$ws1.Range("F7").Value = 200

$ws1.Range("C8").Value = 8100
$ws1.Range("E8").Value = 500
$ws1.Range("F8").Value = 1700

# ---- Sheet 2: "Hospital COVID Census" ----
$ws2.Range("C3").Value = 16

$ws2.Range("C5").Value = 9

$ws2.Range("C6").Value = 104
$ws2.Range("D6").Value = 8

$ws2.Range("C9").Value = 42

$ws2.Range("C10").Value = 19

$ws2.Range("C11").Value = 12
$ws2.Range("D11").Value = 5

$ws2.Range("C12").Value = 21
$ws2.Range("D12").Value = 4

$ws2.Range("C13").Value = 54
$ws2.Range("D13").Value = 11

$ws2.Range("D14").Value = 3

$ws2.Range("C16").Value = 55
$ws2.Range("D16").Value = 12

$ws2.Range("C17").Value = 27
$ws2.Range("D17").Value = 2

$ws2.Range("C18").Value = 68
$ws2.Range("D18").Value = 15

$ws2.Range("C19").Value = 18

$ws2.Range("C20").Value = 31

$ws2.Range("C21").Value = 18
$ws2.Range("D21").Value = 4

$ws2.Range("C22").Value = 30
$ws2.Range("D22").Value = 9

$ws2.Range("D24").Value = 3

$ws2.Range("D26").Value = 1

$ws2.Range("C27").Value = 2

$ws2.Range("C29").Value = 29
$ws2.Range("D29").Value = 5

$ws2.Range("D32").Value = 2

$ws2.Range("D33").Value = 7

$ws2.Range("C34").Value = 16

$ws2.Range("C35").Value = 49

$ws2.Range("C36").Value = 0

$ws2.Range("C37").Value = 35
$ws2.Range("D37").Value = 2

$ws2.Range("C38").Value = 52
$ws2.Range("D38").Value = 10

$ws2.Range("C39").Value = 8
$ws2.Range("D39").Value = 4

$ws2.Range("C42").Value = 95
$ws2.Range("D42").Value = 32

$ws2.Range("C43").Value = 35

$ws2.Range("C44").Value = 26
$ws2.Range("D44").Value = 5

$ws2.Range("C45").Value = 2

$ws2.Range("C46").Value = 17
$ws2.Range("D46").Value = 6

$ws2.Range("C48").Value = 41

$ws2.Range("C49").Value = 37

$ws2.Range("C50").Value = 28
$ws2.Range("D50").Value = 8

$ws2.Range("C55").Value = 52
$ws2.Range("D55").Value = 11

$ws2.Range("C57").Value = 39
$ws2.Range("D57").Value = 7

$ws2.Range("C58").Value = 17

$ws2.Range("C59").Value = 33
$ws2.Range("D59").Value = 7

$ws2.Range("C60").Value = 34

$ws2.Range("C61").Value = 26
$ws2.Range("D61").Value = 3

$ws2.Range("C62").Value = 17

$ws2.Range("C63").Value = 48

$ws2.Range("C64").Value = 16
$ws2.Range("D64").Value = 1

$ws2.Range("C65").Value = 6

$ws2.Range("C66").Value = 29
$ws2.Range("D66").Value = 14

$ws2.Range("C67").Value = 21

$ws2.Range("C68").Value = 39
$ws2.Range("D68").Value = 11

$ws2.Range("C69").Value = 34
$ws2.Range("D69").Value = 19

$ws2.Range("C70").Value = 23
$ws2.Range("D70").Value = 3

# ---- Selection / view updates ----
$ws1.Activate()
$ws1.Range("E14").Select() | Out-Null

$ws2.Activate()
$ws2.Range("J13").Select() | Out-Null
